$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "ramanuj"
$ws.Range("B1").Value = "ramanujasati90@gmail.com"

$ws.Range("A2").Value = "science explorers"
$ws.Range("B2").Value = "scienceexplorers31@gmail.com"

$ws.Range("A3").Value = "lavish"
$ws.Range("B3").Value = "sainilavish398@gmail.com"

$ws.Range("B4").Value = "codieszone@gmail.com"

$ws.Range("A5").Value = "faizan"
$ws.Range("B5").Value = "faizananwar344@gmail.com"

$ws.Range("A6").Value = "lavish2"
$ws.Range("B6").Value = "lavishsaini0110@gmail.com"

$ws.Hyperlinks.Add($ws.Range("B1"), "mailto:ramanujasati90@gmail.com")
$ws.Hyperlinks.Add($ws.Range("B5"), "mailto:faizananwar344@gmail.com")

$ws.Range("A4").Value = "codies zone"

$ws.Range("A4").Select()
